$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 is "profile pic" - set the value in column C to "y" (matches existing C3 pattern)
$ws.Range("C14").Value = "y"

# Update selection to C14 to match the saved view state
$ws.Range("C14").Select()
